# Daily attendance processing - reorder "Recorded By" (column G) values so
# that the "System" token is moved to the front of the comma-separated list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        # Case-sensitive search for an exact "System" token. The -eq/-ne
        # operators are case-insensitive here, so use .Equals() (ordinal,
        # case-sensitive) to distinguish "System" from "system".
        $rest = @()
        $foundSystem = $false
        foreach ($p in $parts) {
            if ((-not $foundSystem) -and $p.Equals("System")) {
                $foundSystem = $true
            } else {
                $rest += $p
            }
        }

        if ($foundSystem) {
            $newParts = @("System") + $rest
            $newVal = $newParts -join ", "
            if (-not $newVal.Equals($val)) {
                $cell.Value2 = $newVal
            }
        }
    }
}
